$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove MuSCs sending-cluster rows (rows 5-7); rows below shift up automatically
$ws.Range("A5:T7").Delete()

# Update recomputed TPM-derived values for remaining rows (2-7)
$ws.Range("G2").Value = 0.7189253333333333
$ws.Range("H2").Value = 2.156776
$ws.Range("I2").Value = 0.375764717934637
$ws.Range("J2").Value = 0.3757647179346369
$ws.Range("M2").Value = 4.108517666666667
$ws.Range("N2").Value = 12.325553
$ws.Range("O2").Value = 0.3056854850030113
$ws.Range("P2").Value = 0.3056854850030113
$ws.Range("Q2").Value = 2.953717433014222
$ws.Range("R2").Value = 26.58345689712799
$ws.Range("S2").Value = 0.1148658200488693
$ws.Range("T2").Value = 0.1148658200488692
$ws.Range("G3").Value = 0.7189253333333333
$ws.Range("H3").Value = 2.156776
$ws.Range("I3").Value = 0.375764717934637
$ws.Range("J3").Value = 0.3757647179346369
$ws.Range("M3").Value = 0.053738
$ws.Range("N3").Value = 0.161214
$ws.Range("O3").Value = 0.003998261155444748
$ws.Range("P3").Value = 0.003998261155444747
$ws.Range("Q3").Value = 0.03863360956266666
$ws.Range("R3").Value = 0.3477024860639999
$ws.Range("S3").Value = 0.001502405475304711
$ws.Range("T3").Value = 0.001502405475304711
$ws.Range("G4").Value = 0.7189253333333333
$ws.Range("H4").Value = 2.156776
$ws.Range("I4").Value = 0.375764717934637
$ws.Range("J4").Value = 0.3757647179346369
$ws.Range("M4").Value = 9.278087
$ws.Range("N4").Value = 27.834261
$ws.Range("O4").Value = 0.6903162538415439
$ws.Range("P4").Value = 0.6903162538415439
$ws.Range("Q4").Value = 6.670251789170666
$ws.Range("R4").Value = 60.03226610253599
$ws.Range("S4").Value = 0.259396492410463
$ws.Range("T4").Value = 0.259396492410463
$ws.Range("G5").Value = 1.194307333333333
$ws.Range("H5").Value = 3.582922
$ws.Range("I5").Value = 0.624235282065363
$ws.Range("J5").Value = 0.624235282065363
$ws.Range("M5").Value = 4.108517666666667
$ws.Range("N5").Value = 12.325553
$ws.Range("O5").Value = 0.3056854850030113
$ws.Range("P5").Value = 0.3056854850030113
$ws.Range("Q5").Value = 4.906832778429556
$ws.Range("R5").Value = 44.16149500586599
$ws.Range("S5").Value = 0.1908196649541421
$ws.Range("T5").Value = 0.1908196649541421
$ws.Range("G6").Value = 1.194307333333333
$ws.Range("H6").Value = 3.582922
$ws.Range("I6").Value = 0.624235282065363
$ws.Range("J6").Value = 0.624235282065363
$ws.Range("M6").Value = 0.053738
$ws.Range("N6").Value = 0.161214
$ws.Range("O6").Value = 0.003998261155444748
$ws.Range("P6").Value = 0.003998261155444747
$ws.Range("Q6").Value = 0.06417968747866667
$ws.Range("R6").Value = 0.5776171873079999
$ws.Range("S6").Value = 0.002495855680140036
$ws.Range("T6").Value = 0.002495855680140036
$ws.Range("G7").Value = 1.194307333333333
$ws.Range("H7").Value = 3.582922
$ws.Range("I7").Value = 0.624235282065363
$ws.Range("J7").Value = 0.624235282065363
$ws.Range("M7").Value = 9.278087
$ws.Range("N7").Value = 27.834261
$ws.Range("O7").Value = 0.6903162538415439
$ws.Range("P7").Value = 0.6903162538415439
$ws.Range("Q7").Value = 11.08088734340467
$ws.Range("R7").Value = 99.72798609064199
$ws.Range("S7").Value = 0.430919761431081
$ws.Range("T7").Value = 0.430919761431081
